# Commit: Wed, Jun 03, 2020  4:04:52 AM
#
# 1) Slide 16 has a 2-column summary table (shape 3) whose table style was
#    changed to a different built-in style ("EC9BE02F-...", a Medium-style
#    accented table look) via the Table Design gallery.
# 2) The presentation's Design/theme was changed from the old "Integral"
#    theme to the built-in "Office Theme" color palette (Design tab ->
#    Variants / theme gallery). That rewrites the 12 theme colors
#    (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) used by the slide master's
#    theme (theme1.xml); fonts and format scheme are already identical
#    between the old and new theme so nothing else changes there.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 ---------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{EC9BE02F-3D67-4A16-8C3A-DD38457F1E29}")

# --- 2. Switch the presentation's theme colors to "Office Theme" -------
# ThemeColorScheme is shared by the whole deck (it lives on the slide
# master's theme part), so it can be read/written from any slide.
$themeSlide = $p.Slides.Item(1)
$colors = $themeSlide.ThemeColorScheme

$colors.Item(1).RGB  = 0        # dk1      000000
$colors.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388  # dk2      44546A
$colors.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501  # accent2  ED7D31
$colors.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colors.Item(8).RGB  = 49407    # accent4  FFC000
$colors.Item(9).RGB  = 12874308 # accent5  4472C4
$colors.Item(10).RGB = 4697456  # accent6  70AD47
$colors.Item(11).RGB = 12673797 # hlink    0563C1
$colors.Item(12).RGB = 7491477  # folHlink 954F72
